# Doing Updates for Financials
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HPJ")

# Row 21 - Non Recurring: oldest year (J, 2011) data no longer available
$ws.Range("J21").Value = "NA"

# Row 44 - Cash And Cash Equivalents: revised figures for 2017 (D) and 2016 (E)
$ws.Range("D44").Value = 49000
$ws.Range("E44").Value = 46400

# Row 45 - Short Term Investments: revised figures for 2017 (D) and 2016 (E)
$ws.Range("D45").Value = 27200
$ws.Range("E45").Value = 19300

# Row 47 - Inventory: revised figure for 2017 (D)
$ws.Range("D47").Value = 9900

# Row 52 - Long Term Investments: revised figure for 2017 (D)
$ws.Range("D52").Value = 7400

# Row 83 - Depreciation: oldest year (J, 2011) data no longer available
$ws.Range("J83").Value = "NA"

# Row 91 - Capital Expenditures: revised figure for 2012 (I)
$ws.Range("I91").Value = -11600

# Row 94 - Total Cash Flows From Investing Activities: oldest year (J, 2011) data no longer available
$ws.Range("J94").Value = "NA"

# Row 100 - Total Cash Flows From Financing Activities: oldest year (J, 2011) data no longer available
$ws.Range("J100").Value = "NA"

# Row 101 - Effect Of Exchange Rate Changes: oldest year (J, 2011) data no longer available
$ws.Range("J101").Value = "NA"
